$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2024-09-05 Thursday" "2024-09-06 Friday"

Replace-Text "36÷5=7, 1" "25÷6=4, 1"
Replace-Text "13÷7=1, 6" "95÷8=11, 7"
Replace-Text "73÷3=24, 1" "59÷7=8, 3"
Replace-Text "76÷2=38, 0" "23÷9=2, 5"
Replace-Text "64÷9=7, 1" "53÷9=5, 8"

Replace-Text "52÷3=17, 1" "85÷8=10, 5"
Replace-Text "76÷3=25, 1" "56÷8=7, 0"
Replace-Text "57÷7=8, 1" "71÷7=10, 1"
Replace-Text "67÷5=13, 2" "10÷4=2, 2"
Replace-Text "26÷7=3, 5" "52÷6=8, 4"

Replace-Text "52÷2=26, 0" "72÷9=8, 0"
Replace-Text "32÷8=4, 0" "98÷9=10, 8"
Replace-Text "80÷4=20, 0" "21÷3=7, 0"
Replace-Text "19÷5=3, 4" "24÷8=3, 0"
Replace-Text "90÷7=12, 6" "31÷8=3, 7"

Replace-Text "32÷2=16, 0" "58÷6=9, 4"
Replace-Text "12÷3=4, 0" "76÷5=15, 1"
Replace-Text "99÷9=11, 0" "94÷5=18, 4"
Replace-Text "91÷2=45, 1" "92÷9=10, 2"
Replace-Text "97÷9=10, 7" "30÷8=3, 6"

Replace-Text "57÷9=6, 3" "35÷7=5, 0"
Replace-Text "70÷6=11, 4" "21÷8=2, 5"
Replace-Text "22÷5=4, 2" "52÷9=5, 7"
Replace-Text "86÷2=43, 0" "36÷9=4, 0"
Replace-Text "45÷6=7, 3" "34÷6=5, 4"
